$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("A7").Value = "Helping with system UI"
$ws.Range("B7").Value = 43747
$ws.Range("B7").NumberFormat = "d-mmm"
$ws.Range("C7").Value = 43748
$ws.Range("C7").NumberFormat = "d-mmm"
$ws.Range("D7").Value = 0.77083333333333337
$ws.Range("D7").NumberFormat = "h:mm AM/PM"
$ws.Range("E7").Value = 0.91666666666666663
$ws.Range("E7").NumberFormat = "h:mm AM/PM"

# Row 8
$ws.Range("A8").Value = "Helping with Documentation of the System"
$ws.Range("B8").Value = 43753
$ws.Range("B8").NumberFormat = "d-mmm"
$ws.Range("C8").Value = 43755
$ws.Range("C8").NumberFormat = "d-mmm"
$ws.Range("D8").Value = 0.79166666666666663
$ws.Range("D8").NumberFormat = "h:mm AM/PM"
$ws.Range("E8").Value = 0.91666666666666663
$ws.Range("E8").NumberFormat = "h:mm AM/PM"

# Row 9
$ws.Range("A9").Value = "Helping with further development of system UI"
$ws.Range("B9").Value = 43765
$ws.Range("B9").NumberFormat = "d-mmm"
$ws.Range("C9").Value = 43766
$ws.Range("C9").NumberFormat = "d-mmm"
$ws.Range("D9").Value = 0.625
$ws.Range("D9").NumberFormat = "h:mm AM/PM"
$ws.Range("E9").Value = 0.91666666666666663
$ws.Range("E9").NumberFormat = "h:mm AM/PM"

# Update selection to match the final state recorded in the workbook
$ws.Range("A16").Select()
